# Evin Lewis.xlsx edit: rename sheet, prepend "matchNo" column, and add
# four more match rows of scraped batting data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet tab.
$ws.Name = "Evin Lewis"

# 2) Clear everything so we can lay the table out fresh (new column layout).
$ws.Cells.Clear()

# 3) Header row.
$headers = @("matchNo","teamName","batterName","states","runs","balls","fours","sixes","sr","opponentTeamName","venue","date","result")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# 4) Data rows (all values are text in the source data, including the
#    numeric-looking ones, so every cell is a string).
$rows = @(
    @("47th","Rajasthan Royals","Evin Lewis","c Hazlewood b Thakur","27","12","2","2","225.00","Chennai Super Kings","Abu Dhabi","October 02","Royals won by 7 wickets (with 15 balls remaining)"),
    @("51st","Rajasthan Royals","Evin Lewis","lbw b Bumrah","24","19","3","1","126.31","Mumbai Indians","Sharjah","October 05","Mumbai won by 8 wickets (with 70 balls remaining)"),
    @("43rd","Rajasthan Royals","Evin Lewis","c †Bharat b Garton","58","37","5","3","156.75","Royal Challengers Bangalore","Dubai (DSC)","September 29","RCB won by 7 wickets (with 17 balls remaining)"),
    @("32nd","Rajasthan Royals","Evin Lewis","c Agarwal b Arshdeep Singh","36","21","7","1","171.42","Punjab Kings","Dubai (DSC)","September 21","Royals won by 2 runs"),
    @("40th","Rajasthan Royals","Evin Lewis","c Abdul Samad b Kumar","6","4","1","0","150.00","Sunrisers Hyderabad","Dubai (DSC)","September 27","Sunrisers won by 7 wickets (with 9 balls remaining)")
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowData = $rows[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $cell = $ws.Cells.Item($r + 2, $c + 1)
        $text = $rowData[$c]
        $looksNumeric = $text -match '^-?[0-9]+(\.[0-9]+)?$'
        if ($looksNumeric) {
            # Force text storage (leading apostrophe) but strip the
            # quote-prefix style afterwards so the cell keeps the sheet's
            # default styling, matching the source file.
            $cell.Value = "'" + $text
            $cell.Style = "Normal"
        } else {
            $cell.Value = $text
        }
    }
}

Write-Output "Evin Lewis sheet rebuilt"
